# Final code changes on 22nd feb 2022
#
# Job# value in B2 is updated from "72004073" to "32265125". The value is a
# digit-only string that must stay stored as text (as it already was, along
# with B3/B4 in the same column), not be auto-converted to a number.
#
# Entering a leading apostrophe forces Excel to treat the numeric-looking
# text as a string, but it also stamps the cell with a "quote prefix" number
# format. The sibling cells (B3/B4) hold the same kind of text without any
# such formatting, so we copy their (default) formatting back onto B2 after
# the value is set, restoring its original unformatted appearance while
# keeping the new value as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "'32265125"

$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
